$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (week of 2021-09-22, serial 44461) is
# inserted into the "Perejil" data block. It becomes the new row 194,
# pushing the former rows 194-211 down to 195-212 (dimension grows from
# A1:R211 to A1:R212).
$ws.Rows.Item(194).Insert()

$ws.Cells.Item(194, 1).Value = 9
$ws.Cells.Item(194, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(194, 3).Value = "Metropolitana"
$ws.Cells.Item(194, 4).Value = 44461
$ws.Cells.Item(194, 5).Value = 13
$ws.Cells.Item(194, 6).Value = 100112044
$ws.Cells.Item(194, 7).Value = "Perejil"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 79
$ws.Cells.Item(194, 11).Value = 13000
$ws.Cells.Item(194, 12).Value = 14000
$ws.Cells.Item(194, 13).Value = 13494
$ws.Cells.Item(194, 14).Value = "$/docena de atados"
$ws.Cells.Item(194, 15).Value = "Región Metropolitana"
$ws.Cells.Item(194, 16).Value = 4498
$ws.Cells.Item(194, 17).Value = 3
$ws.Cells.Item(194, 18).Value = "Hortaliza"
